# Palred Technologies - Quarterly sheet: insert a new "Exceptional items" column
# before the existing "P/l before tax" column (column L), shifting columns
# L..T one position to the right (M..U) for the data rows, and appending a
# duplicated header cell in U1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Map of row -> new numeric value to place in the freshly inserted column L
# (only rows that actually report a non-zero "Exceptional items" figure get a
# value here; all other rows are left blank).
$overrides = @{
    3  = 166.58
    7  = -0.98
    17 = -0.78
    23 = -30.91
    24 = -10
    27 = -6.87
    29 = -1.42
    34 = -0.49
    36 = -3.33
    37 = 4.41
    47 = -2.01
}

# Shift columns L(12) through T(20) one column to the right, into M(13)..U(21),
# for every data row (2..47). Work right-to-left so we never overwrite a
# value before it has been read.
for ($row = 2; $row -le 47; $row++) {
    for ($col = 20; $col -ge 12; $col--) {
        $srcCell = $ws.Cells.Item($row, $col)
        $dstCell = $ws.Cells.Item($row, $col + 1)
        $val = $srcCell.Value2
        if ($val -eq $null) {
            $dstCell.ClearContents()
        } else {
            $dstCell.Value2 = $val
        }
    }

    # Populate the newly freed column L (12).
    $lCell = $ws.Cells.Item($row, 12)
    if ($row -eq 2) {
        $lCell.Value2 = "Exceptional Items"
    } elseif ($overrides.ContainsKey($row)) {
        $lCell.Value2 = $overrides[$row]
    } else {
        $lCell.ClearContents()
    }
}

# Row 1 (column headers) is left untouched except for a new trailing cell
# U1, which simply duplicates the last existing header (T1), formatting
# included.
$t1 = $ws.Cells.Item(1, 20)
$u1 = $ws.Cells.Item(1, 21)
$t1.Copy($u1)
